# Laborator 14.03.2024 - Fake Photoshop. Gray scale, Complementary, Contrast, Blur.
# Mark attendance for "săpt. 3" (column E) for the students that were present,
# add the three new students that joined this lab, and update the
# view/selection to where the professor left off (row 16 / cell E26).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Mark column E ("săpt. 3") attendance as TRUE for the students present ---
$presentRows = @(4, 7, 8, 13, 18, 19, 20, 22, 25, 26, 29, 32, 33, 34)
foreach ($r in $presentRows) {
    $ws.Cells.Item($r, 5).Value = $true   # column E = 5
}

# --- Three new students joined; add their names and mark them present too ---
$ws.Range("B32").Value = "Raluca Veres"
$ws.Range("B33").Value = "Gabriela Maghear"
$ws.Range("B34").Value = "Sebastian Pop"

# --- Scroll the sheet so row 16 is at the top, and leave the selection on E26 ---
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E26").Select()
